$d = $word.ActiveDocument

# Test: replace Globality -> Module using Find/Replace
$d.Content.Find.Execute("Globality", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Module", 2)
